$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. After the "Laid before Parliament" paragraph, insert two new date-block
#    paragraphs: "LaidDraft" (Unknown looks like a date) and "Sifted"
#    (Sift requirements satisfied).
# ---------------------------------------------------------------------------
$laid = $d.Paragraphs.Item(6)
if ($laid.Range.Text -notlike "Laid before Parliament*") {
    throw "Paragraph 6 is not the 'Laid before Parliament' paragraph: $($laid.Range.Text)"
}

[void]$laid.Range.InsertParagraphAfter()
$laidDraft = $d.Paragraphs.Item(7)
$laidDraftXml = "<w:p $wns><w:pPr><w:pStyle w:val='LaidDraft'/></w:pPr>" +
    "<w:r><w:t>Unknown looks like a date</w:t></w:r>" +
    "<w:r><w:tab/></w:r>" +
    "<w:r><w:t>2</w:t></w:r>" +
    "<w:r><w:t>4</w:t></w:r>" +
    "<w:r><w:t>th September 2021</w:t></w:r>" +
    "</w:p>"
[void]$laidDraft.Range.InsertXML($laidDraftXml)

$laidDraft = $d.Paragraphs.Item(7)
[void]$laidDraft.Range.InsertParagraphAfter()
$sifted = $d.Paragraphs.Item(8)
$siftedXml = "<w:p $wns><w:pPr><w:pStyle w:val='Sifted'/></w:pPr>" +
    "<w:r><w:t>Sift requirements satisfied</w:t></w:r>" +
    "<w:r><w:tab/></w:r>" +
    "<w:r><w:t>30</w:t></w:r>" +
    "<w:r><w:t>th September 2021</w:t></w:r>" +
    "</w:p>"
[void]$sifted.Range.InsertXML($siftedXml)

# ---------------------------------------------------------------------------
# 2. Shrink the "Coming into force" paragraph down to just "Coming into
#    force" (split as "forc"+"e"), dropping the old tab/dash placeholders
#    and date. Follow it with two new "ComingC" paragraphs giving the real
#    commencement dates.
# ---------------------------------------------------------------------------
$coming = $d.Paragraphs.Item(9)
if ($coming.Range.Text -notlike "Coming into force*") {
    throw "Paragraph 9 is not the 'Coming into force' paragraph: $($coming.Range.Text)"
}

$comingXml = "<w:p $wns><w:pPr><w:pStyle w:val='Coming'/></w:pPr>" +
    "<w:r><w:t>Coming</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t>into</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t>forc</w:t></w:r>" +
    "<w:r><w:t>e</w:t></w:r>" +
    "</w:p>"
[void]$coming.Range.InsertXML($comingXml)

$coming = $d.Paragraphs.Item(9)
[void]$coming.Range.InsertParagraphAfter()
$comingC1 = $d.Paragraphs.Item(10)
$comingC1Xml = "<w:p $wns><w:pPr><w:pStyle w:val='ComingC'/></w:pPr>" +
    "<w:r><w:t>For the purpose of</w:t></w:r>" +
    "<w:r><w:tab/></w:r>" +
    "<w:r><w:t>20th September 2021</w:t></w:r>" +
    "</w:p>"
[void]$comingC1.Range.InsertXML($comingC1Xml)

$comingC1 = $d.Paragraphs.Item(10)
[void]$comingC1.Range.InsertParagraphAfter()
$comingC2 = $d.Paragraphs.Item(11)
$comingC2Xml = "<w:p $wns><w:pPr><w:pStyle w:val='ComingC'/></w:pPr>" +
    "<w:r><w:t>For the</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> other</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> purpose of</w:t></w:r>" +
    "<w:r><w:tab/></w:r>" +
    "<w:r><w:t>2</w:t></w:r>" +
    "<w:r><w:t>1st</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> September 2021</w:t></w:r>" +
    "</w:p>"
[void]$comingC2.Range.InsertXML($comingC2Xml)

# ---------------------------------------------------------------------------
# 3. Move the "lastRenderedPageBreak" marker: it now sits before "In" in the
#    "In regulation 9(3) (renewal of registration documents)..." paragraph,
#    and is removed from the "The amendment in the Schedule has effect."
#    paragraph (whose bookmark must be preserved).
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "In regulation 9(3)*") {
        $found = $true
        $xml = "<w:p $wns><w:pPr><w:pStyle w:val='N2'/></w:pPr>" +
            "<w:r><w:lastRenderedPageBreak/><w:t>In</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>regulation</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>9(3)</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>(renewal</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>of</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>registration</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>documents)</w:t></w:r>" +
            "<w:r><w:t>,</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>for</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>&#8220;</w:t></w:r>" +
            "<w:r><w:t>one</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>calendar</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>year</w:t></w:r>" +
            "<w:r><w:t>&#8221;</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>substitute</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>&#8220;six</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>months&#8221;.</w:t></w:r>" +
            "</w:p>"
        $pp.Range.InsertXML($xml)
        break
    }
}
if (-not $found) { throw "Could not find the 'In regulation 9(3)' paragraph" }

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "The amendment in the Schedule has effect*") {
        $found = $true
        $xml = "<w:p $wns><w:pPr><w:pStyle w:val='N2'/></w:pPr>" +
            "<w:bookmarkStart w:id='0' w:name='_Ref70071194'/>" +
            "<w:r><w:t>The</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>amendment</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>in</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>the</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>Schedule</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>has</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
            "<w:r><w:t>effect.</w:t></w:r>" +
            "<w:bookmarkEnd w:id='0'/>" +
            "</w:p>"
        $pp.Range.InsertXML($xml)
        break
    }
}
if (-not $found) { throw "Could not find the 'The amendment in the Schedule has effect.' paragraph" }

Write-Host "Done"
